$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert the new work-entry data for row 28 (a new day worked: 18.9.2020, 14:00-20:00)
# Copy the number formatting from the row above (row 27) so the new cells
# render as date / time values, just like the rest of the table.
$ws.Range("A27:C27").Copy()
$ws.Range("A28:C28").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A28").Value = 44092
$ws.Range("B28").Value = 0.58333333333333337
$ws.Range("C28").Value = 0.83333333333333337

# Move the active selection to C29 like in the edited workbook
$ws.Range("C29").Select()
